$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (row 2 through row 463) holds the "Förändrad" (last changed) date
# serial value 46081 (2026-02-28). Bump it by one day to 46082 (2026-03-01)
# for every data row, matching the commit's automatic daily update.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162
if ($lastRow -lt 2) { $lastRow = 463 }

$range = $ws.Range("C2:C$lastRow")
$range.Value = 46082
